$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples")

# The sheet currently has header row 1, data rows 2-3, and a stray cell at G6.
# We need 6 full data rows (2-7): B1/B2/B3 wells for Positive sequence and
# C1/C2/C3 wells for Negative sequence, each with Tune Mix / IPA Blank / Sample
# descriptions (no more separate calibrant run row).

# Remove the stray leftover row first.
$ws.Rows.Item(6).Delete()

# Insert four fresh rows (4,5,6,7), copying row 2's formatting so that every
# cell keeps the same style used across the table.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(4).Insert(-4121)
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(5).Insert(-4121)
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(6).Insert(-4121)
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(7).Insert(-4121)

# Row 2: B1 - Tune Mix - Positive
$ws.Range("A2").Value = "B1"
$ws.Range("B2").NumberFormat = "General"
$ws.Range("B2").Value = "Tune Mix"
$ws.Range("C2").Value = "Positive"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "TUNE"
$ws.Range("G2").Value = "2023-03-02_dodd_4bit_POS.m"
$ws.Range("H2").Value = "P384"
$ws.Range("I2").Value = "C"

# Row 3: B2 - IPA Blank - Positive
$ws.Range("A3").Value = "B2"
$ws.Range("B3").Value = "IPA Blank"
$ws.Range("C3").Value = "Positive"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "BLANK"
$ws.Range("G3").Value = "2023-03-02_dodd_4bit_POS.m"
$ws.Range("H3").Value = "P384"
$ws.Range("I3").Value = "C"

# Row 4: B3 - Sample - Positive
$ws.Range("A4").Value = "B3"
$ws.Range("B4").Value = "Sample"
$ws.Range("C4").Value = "Positive"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "SAMPLE"
$ws.Range("G4").Value = "2023-03-02_dodd_4bit_POS.m"
$ws.Range("H4").Value = "P384"
$ws.Range("I4").Value = "C"

# Row 5: C1 - Tune Mix - Negative
$ws.Range("A5").Value = "C1"
$ws.Range("B5").NumberFormat = "General"
$ws.Range("B5").Value = "Tune Mix"
$ws.Range("C5").Value = "Negative"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "TUNE"
$ws.Range("G5").Value = "2023-03-02_dodd_4bit_NEG.m"
$ws.Range("H5").Value = "P384"
$ws.Range("I5").Value = "C"

# Row 6: C2 - IPA Blank - Negative
$ws.Range("A6").Value = "C2"
$ws.Range("B6").Value = "IPA Blank"
$ws.Range("C6").Value = "Negative"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "BLANK"
$ws.Range("G6").Value = "2023-03-02_dodd_4bit_NEG.m"
$ws.Range("H6").Value = "P384"
$ws.Range("I6").Value = "C"

# Row 7: C3 - Sample - Negative
$ws.Range("A7").Value = "C3"
$ws.Range("B7").Value = "Sample"
$ws.Range("C7").Value = "Negative"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "SAMPLE"
$ws.Range("G7").Value = "2023-03-02_dodd_4bit_NEG.m"
$ws.Range("H7").Value = "P384"
$ws.Range("I7").Value = "C"

# Update selection to reflect the new active cell below the table.
$ws.Range("G8").Select()
